# Bewertungsbogen weights -> DG modell
# Resets the "Note" (grade) inputs in the detail sheet and a few Bonus
# inputs back to blank / zero, tweaks one weighting factor, and fixes a
# handful of typos (double spaces, spaced slashes, "PDF From" -> "PDF-Form")
# in the descriptive texts.

$wb = $excel.ActiveWorkbook

$wsGesamt = $wb.Worksheets.Item("1. Gesamtbewertung (Text)")
$wsDetail = $wb.Worksheets.Item("2. Detailbewertung (Excel)")
$wsKomm   = $wb.Worksheets.Item("3. Kommunikation der Bewertung")

# ---------------------------------------------------------------------
# 1) Text / typo fixes on sheet "2. Detailbewertung (Excel)"
#    (applied in this exact order so the shared-string table ends up in
#    the same order as the authored edit)
# ---------------------------------------------------------------------

$wsDetail.Range("B14").Value = " Analyse von Ergebnissen"
$wsDetail.Range("B18").Value = "DOKUMENTATION, WISSENSTRANSFER"
$wsDetail.Range("B16").Value = "Selbstständigkeit/Betreuungsintensität"

# B7: rich text "Lösungskonzept/Strategie" (title run) + italic explanation run
$b7Title = "Lösungskonzept/Strategie" + [char]10
$b7Body  = "Je nach Komplexität der Aufgabenstellung soll die Gewichtung zwischen 0.2 (einfach) und 1 (komplex) liegen."
$b7Cell = $wsDetail.Range("B7")
$b7Cell.Value = $b7Title + $b7Body
$b7Run2 = $b7Cell.Characters($b7Title.Length + 1, $b7Body.Length)
$b7Run2.Font.ColorIndex = -4105
$b7Run2.Font.Bold = $false
$b7Run2.Font.Italic = $true
$b7Run2.Font.Size = 10
$b7Run2.Font.Name = "Arial"

# B3: rich text "Bemerkungen:" (bold run) + normal explanation run (PDF-Form fix)
$b3Title = "Bemerkungen:"
$b3Body  = " Dieser Bewertungsbogen wird von der betreunden Person ausgefüllt. Bei zwei betreuenden Personen wird er von beiden unabhängig ausgefüllt und danach abgeglichen. Wo möglich und sinnvoll wird ein Kommentar zu jeder Bewertung verfasst. Die Studierenden erhalten in jedem Fall die Würdigung in Papierform. Falls erwünscht wird auch der Bewertungsbogen in PDF-Form abgegeben. Nach der Projektarbeit 5 muss dieser Bewertungsbogen zwingen mit den Studierenden besprochen und auf mögliches Verbesserungspotential für die kommende Projektarbeit 6 hingewiesen werden. Nach Abschluss der Projektarbeit 6 wird der Bewertungsbogen auf Wunsch der Studierenden mit diesen besprochen."
$b3Cell = $wsDetail.Range("B3")
$b3Cell.Value = $b3Title + $b3Body
$b3Run1 = $b3Cell.Characters(1, $b3Title.Length)
$b3Run1.Font.ColorIndex = -4105
$b3Run1.Font.Bold = $true
$b3Run1.Font.Italic = $false
$b3Run1.Font.Size = 14
$b3Run1.Font.Name = "Arial"
$b3Run2 = $b3Cell.Characters($b3Title.Length + 1, $b3Body.Length)
$b3Run2.Font.ColorIndex = -4105
$b3Run2.Font.Bold = $false
$b3Run2.Font.Italic = $false
$b3Run2.Font.Size = 14
$b3Run2.Font.Name = "Arial"

# ---------------------------------------------------------------------
# 2) Weighting / bonus value changes on sheet "2. Detailbewertung (Excel)"
# ---------------------------------------------------------------------

$wsDetail.Range("C14").Value = 1

$wsDetail.Range("C30").Value = 0
$wsDetail.Range("C31").Value = 0
$wsDetail.Range("C32").Value = 0

# ---------------------------------------------------------------------
# 3) Clear out all the "Note" (grade) entries -> blank cells, keeping
#    their existing number formatting / style.
# ---------------------------------------------------------------------

$wsDetail.Range("D7:D9").ClearContents()
$wsDetail.Range("D12:D16").ClearContents()
$wsDetail.Range("D19:D21").ClearContents()
$wsDetail.Range("D24:D26").ClearContents()

# ---------------------------------------------------------------------
# 4) Selection / active-cell bookkeeping to match the saved view state.
# ---------------------------------------------------------------------

$wsGesamt.Activate()
$wsGesamt.Range("B8:C8").Select()

$wsDetail.Activate()
$wsDetail.Range("D7").Select()

$wsKomm.Activate()
